$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data updates (October 11 "tarikh" entry) ---
$ws.Range("B10").Value = 1623
$ws.Range("C10").Value = 977
$ws.Range("B11").Value = 2015

$ws.Range("F14").Value = 100
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 2
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 2
$ws.Range("P14").Value = 2
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = 2
$ws.Range("S14").Value = 2
$ws.Range("T14").Value = 3

# --- Highlight the "total" / "Prev Due to meal" header cells ---
$c2 = $ws.Range("C2")
$c2.Font.Bold = $true
$c2.Font.Size = 12
$c2.Font.Color = 255
$c2.HorizontalAlignment = -4108

$d2 = $ws.Range("D2")
$d2.Font.Bold = $true
$d2.Font.Size = 12
$d2.Font.Color = 5287936
$d2.HorizontalAlignment = -4108

# --- Move the selection / active cell ---
$ws.Range("V2").Select() | Out-Null
